$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 45882
$ws.Range("J3").Value = 45882
$ws.Range("L3").Value = 45882
$ws.Range("N3").Value = -46110

$ws.Range("H39").Value = 258.23077
$ws.Range("I39").Value = 65.7
$ws.Range("K39").Value = 197.1
$ws.Range("M39").Value = 98.89999999999998

$ws.Range("H97").Value = 47166
$ws.Range("J97").Value = 47166
$ws.Range("L97").Value = 141498
$ws.Range("N97").Value = -142490

$ws.Range("H98").Value = 1224.52
$ws.Range("I98").Value = 1161.174
$ws.Range("J98").Value = 1953
$ws.Range("K98").Value = 1161.174
$ws.Range("L98").Value = 1953
$ws.Range("M98").Value = 336.826
$ws.Range("N98").Value = -4949

$ws.Range("H102").Value = 45882
$ws.Range("J102").Value = 45882
$ws.Range("L102").Value = 45882
$ws.Range("N102").Value = -52372

$ws.Range("H113").Value = 57411.332
$ws.Range("J113").Value = 10903.454
$ws.Range("L113").Value = 10903.454
$ws.Range("N113").Value = -17411.454

$ws.Range("H122").Value = 1224.52
$ws.Range("I122").Value = 1161.174
$ws.Range("J122").Value = 1953
$ws.Range("K122").Value = 3483.522
$ws.Range("L122").Value = 5859
$ws.Range("M122").Value = -1033.522
$ws.Range("N122").Value = -10759

$ws.Range("H132").Value = 1392.4916
$ws.Range("I132").Value = 1107.3877
$ws.Range("K132").Value = 3322.1631
$ws.Range("M132").Value = -792.1630999999998

$ws.Range("H135").Value = 1204.7059
$ws.Range("I135").Value = 1030
$ws.Range("K135").Value = 9270
$ws.Range("M135").Value = -6735

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4531.241
$ws.Range("J32").Value = 11883.714
$ws.Range("L32").Value = 11883.714
$ws.Range("N32").Value = -12457.714

$ws.Range("H69").Value = 500224.5
$ws.Range("J69").Value = 500224.5
$ws.Range("L69").Value = 500224.5
$ws.Range("N69").Value = -501722.5

$ws.Range("H72").Value = 500224.5
$ws.Range("J72").Value = 500224.5
$ws.Range("L72").Value = 1500673.5
$ws.Range("N72").Value = -1508161.5

$ws.Range("H97").Value = 591.0323
$ws.Range("I97").Value = 664.5769
$ws.Range("J97").Value = 208.6
$ws.Range("K97").Value = 664.5769
$ws.Range("L97").Value = 208.6
$ws.Range("M97").Value = -168.5769
$ws.Range("N97").Value = -1200.6

$ws.Range("H101").Value = 25000
$ws.Range("J101").Value = 25000
$ws.Range("L101").Value = 25000
$ws.Range("N101").Value = -31490

$ws.Range("H122").Value = 2534.5
$ws.Range("I122").Value = 2336.889
$ws.Range("J122").Value = 3423.75
$ws.Range("K122").Value = 7010.667
$ws.Range("L122").Value = 10271.25
$ws.Range("M122").Value = -4560.667
$ws.Range("N122").Value = -15171.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -21232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 769.6667
$ws.Range("I14").Value = 769.6667
$ws.Range("K14").Value = 769.6667
$ws.Range("M14").Value = -599.6667

$ws.Range("H28").Value = 20011.834
$ws.Range("J28").Value = 19414.4
$ws.Range("L28").Value = 19414.4
$ws.Range("N28").Value = -19904.4

$ws.Range("H43").Value = 13130
$ws.Range("J43").Value = 13130
$ws.Range("L43").Value = 13130
$ws.Range("N43").Value = -13498

$ws.Range("H99").Value = 10833
$ws.Range("J99").Value = 10833
$ws.Range("L99").Value = 10833
$ws.Range("N99").Value = -13829

$ws.Range("H101").Value = 13130
$ws.Range("J101").Value = 13130
$ws.Range("L101").Value = 13130
$ws.Range("N101").Value = -19620

$ws.Range("H126").Value = 10833
$ws.Range("J126").Value = 10833
$ws.Range("L126").Value = 32499
$ws.Range("N126").Value = -37439

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 5195.6924
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 5587
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 16761
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -17099

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H27").Value = 5195.6924
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 5587
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 16761
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -16965

$ws.Range("H38").Value = 83.5
$ws.Range("I38").Value = 77.833336
$ws.Range("J38").Value = 100.5
$ws.Range("K38").Value = 233.500008
$ws.Range("L38").Value = 301.5
$ws.Range("M38").Value = 113.499992
$ws.Range("N38").Value = -995.5

$ws.Range("H92").Value = 1673.3334
$ws.Range("I92").Value = 320
$ws.Range("K92").Value = 960
$ws.Range("M92").Value = 288

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 2497.2666
$ws.Range("I122").Value = 1399.1428
$ws.Range("K122").Value = 4197.428400000001
$ws.Range("M122").Value = -1747.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1550.3
$ws.Range("I22").Value = 1555.3334
$ws.Range("J22").Value = 1549.4117
$ws.Range("K22").Value = 1555.3334
$ws.Range("L22").Value = 1549.4117
$ws.Range("M22").Value = -1260.3334
$ws.Range("N22").Value = -2139.4117

$ws.Range("H27").Value = 1550.3
$ws.Range("I27").Value = 1555.3334
$ws.Range("J27").Value = 1549.4117
$ws.Range("K27").Value = 1555.3334
$ws.Range("L27").Value = 1549.4117
$ws.Range("M27").Value = -1448.3334
$ws.Range("N27").Value = -1763.4117

$ws.Range("H46").Value = 1979.8182
$ws.Range("I46").Value = 1769.8572
$ws.Range("J46").Value = 2347.25
$ws.Range("K46").Value = 1769.8572
$ws.Range("L46").Value = 2347.25
$ws.Range("M46").Value = -1581.8572
$ws.Range("N46").Value = -2723.25

$ws.Range("H55").Value = 784.2857
$ws.Range("I55").Value = 400
$ws.Range("K55").Value = 400
$ws.Range("M55").Value = -227

$ws.Range("H61").Value = 3892
$ws.Range("I61").Value = 3833.6843
$ws.Range("K61").Value = 3833.6843
$ws.Range("M61").Value = -3631.6843

$ws.Range("H113").Value = 3892
$ws.Range("I113").Value = 3833.6843
$ws.Range("K113").Value = 3833.6843
$ws.Range("M113").Value = -1663.6843

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 347565.38
$ws.Range("I14").Value = 427478.8
$ws.Range("K14").Value = 427478.8
$ws.Range("M14").Value = -427310.8

$ws.Range("H15").Value = 55003
$ws.Range("J15").Value = 60000
$ws.Range("L15").Value = 60000
$ws.Range("N15").Value = -60576

$ws.Range("H19").Value = 4000
$ws.Range("I19").Value = 4000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -3826
$ws.Range("N19").ClearContents()

$ws.Range("H62").Value = 14070.333
$ws.Range("J62").Value = 10765.346
$ws.Range("L62").Value = 10765.346
$ws.Range("N62").Value = -12013.346

$ws.Range("H65").Value = 14070.333
$ws.Range("J65").Value = 10765.346
$ws.Range("L65").Value = 53826.73
$ws.Range("N65").Value = -60066.73

$ws.Range("H95").Value = 60039.715
$ws.Range("J95").Value = 60039.715
$ws.Range("L95").Value = 60039.715
$ws.Range("N95").Value = -65531.715

$ws.Range("H122").Value = 2495.0466
$ws.Range("I122").Value = 2310.1614
$ws.Range("K122").Value = 6930.4842
$ws.Range("M122").Value = -4480.4842

$ws.Range("H132").Value = 1522.5428
$ws.Range("I132").Value = 1499.8334
$ws.Range("K132").Value = 4499.5002
$ws.Range("M132").Value = -1969.5002

$ws.Range("H139").Value = 103503
$ws.Range("J139").Value = 103503
$ws.Range("L139").Value = 103503
$ws.Range("N139").Value = -113783
